$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 5 new attendee rows (85-89) to the WG attendees sheet, matching the
# "Minutes for 10-04-2024" update: two more R Consortium folks, a new Sanofi
# attendee, a new J&J attendee and a new Novo Nordisk attendee.
# ---------------------------------------------------------------------------

# Adding a hyperlink via the COM `Hyperlinks.Add` call stamps the cell with
# a brand-new "Hyperlink" named style (underline + theme color), which isn't
# what this sheet uses -- every existing affiliation_url cell is plain Arial
# 10pt in blue (FF0000FF), using the workbook's existing "s=2" cell style.
# So: add the hyperlink first, then paste the formatting from an existing
# affiliation_url cell (C84) back on top to restore the sheet's look.
function Set-AffiliationUrl($row, $url, $display) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $url
    $ws.Hyperlinks.Add($cell, $url, "", "", $display)
    $ws.Range("C84").Copy()
    $cell.PasteSpecial(-4122)
}

# Row 85: Terry Christiani / R Consortium
$ws.Cells.Item(85, 1).Value = "Terry Christiani"
$ws.Cells.Item(85, 2).Value = "R Consortium"
Set-AffiliationUrl 85 "https://www.r-consortium.org" "https://www.r-consortium.org"

# Row 86: Amanda Martin / R Consortium
$ws.Cells.Item(86, 1).Value = "Amanda Martin"
$ws.Cells.Item(86, 2).Value = "R Consortium"
Set-AffiliationUrl 86 "https://www.r-consortium.org" "https://www.r-consortium.org"

# Row 87: Andre Couturier / Sanofi (affiliation url has a trailing space run
# with its own, unformatted, run -- matches the source workbook's pattern)
$ws.Cells.Item(87, 1).Value = "Andre Couturier"
$ws.Cells.Item(87, 2).Value = "Sanofi"
$ws.Cells.Item(87, 3).Value = "https://www.sanofi.com/ "
$ws.Hyperlinks.Add($ws.Cells.Item(87, 3), "https://www.sanofi.com/", "", "", "https://www.sanofi.com/")
$urlPart = $ws.Cells.Item(87, 3).Characters(1, 23)
$urlPart.Font.Name = "Arial"
$urlPart.Font.Size = 10
$urlPart.Font.Color = 16711680
$spacePart = $ws.Cells.Item(87, 3).Characters(24, 1)
$spacePart.Font.Name = "Arial"
$spacePart.Font.Size = 10
$spacePart.Font.ColorIndex = -4105

# Row 88: Nicholas Masel / Johnson & Johnson
$ws.Cells.Item(88, 1).Value = "Nicholas Masel"
$ws.Cells.Item(88, 2).Value = "Johnson & Johnson"
Set-AffiliationUrl 88 "https://www.jnj.com" "https://www.jnj.com"

# Row 89: Lovemore Gakava / Novo Nordisk
$ws.Cells.Item(89, 1).Value = "Lovemore Gakava"
$ws.Cells.Item(89, 2).Value = "Novo Nordisk"
Set-AffiliationUrl 89 "https://www.novonordisk.com/ " "https://www.novonordisk.com/ "

# Drop the unused "Hyperlink" named cell style that `Hyperlinks.Add` creates
# behind the scenes so the style table matches the rest of the sheet.
$wb.Styles("Hyperlink").Delete()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Misc sheet view bookkeeping that accompanied the data edit upstream.
# ---------------------------------------------------------------------------
$ws.Range("C1:C89").ColumnWidth = 36.3
$ws.Range("A70").Select()
